$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (Invoice Number, Date of Invoice, Payment Date, Amount)
$data = @(
    @(6510848763, "2016.11.07", "2017.01.11", 266915),
    @(6510848803, "2016.11.07", "2017.01.09", 245895),
    @(6510848932, "2016.11.07", "2017.01.09", 1018340),
    @(6510849001, "2016.11.08", "2017.01.11", 964590),
    @(6510849079, "2016.11.08", "2017.01.11", 248455),
    @(6510849149, "2016.11.08", "2017.01.11", 1074500),
    @(6510849460, "2016.11.07", "2017.01.11", 36895),
    @(6510849534, "2016.11.08", "2017.01.11", 26085)
)

# Clear out the old data rows (rows 2 through 16) first
$ws.Range("B2:E16").Clear()

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $r = $r + 1
}
